$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Förändrad" (column C) date bumped from 2023-09-03 (45172) to 2023-09-06 (45175)
#    for every existing data row (2..54).
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 3).Value = 45175
}

# 2) Row 54 gains an explicit 15pt custom row height.
$ws.Rows.Item(54).RowHeight = 15

# 3) Three new data rows are appended (55-57); rows 55 & 56 carry an explicit
#    15pt custom row height, row 57 keeps the sheet's default height.
$newRows = @(
    @{ Row = 55; A = "A 41272-2023"; B = 45174; C = 45175; G = 0.9; SetHeight = $true },
    @{ Row = 56; A = "A 41263-2023"; B = 45174; C = 45175; G = 4.7; SetHeight = $true },
    @{ Row = 57; A = "A 41276-2023"; B = 45174; C = 45175; G = 0.3; SetHeight = $false }
)

foreach ($row in $newRows) {
    $r = $row.Row

    $ws.Cells.Item($r, 1).Value = $row.A

    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 2).Value = $row.B

    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 3).Value = $row.C

    $ws.Cells.Item($r, 4).Value = "UPPSALA LÄN"
    $ws.Cells.Item($r, 5).Value = "ÄLVKARLEBY"
    $ws.Cells.Item($r, 6).Value = "Bergvik skog väst AB"

    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = 0
    $ws.Cells.Item($r, 10).Value = 0
    $ws.Cells.Item($r, 11).Value = 0
    $ws.Cells.Item($r, 12).Value = 0
    $ws.Cells.Item($r, 13).Value = 0
    $ws.Cells.Item($r, 14).Value = 0
    $ws.Cells.Item($r, 15).Value = 0
    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0

    # Column R stays an empty, wrap-text styled cell (same style as the rest
    # of the table's R column).
    $ws.Cells.Item($r, 18).WrapText = $true

    if ($row.SetHeight) {
        $ws.Rows.Item($r).RowHeight = 15
    }
}
